$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.61
$ws.Range("Q2").Value = 2.01
$ws.Range("R2").Value = 1.71
$ws.Range("G3").Value = 1.81
$ws.Range("Q3").Value = 2.55
$ws.Range("R3").Value = 1.44
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 41
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 5.5
$ws.Range("J4").Value = 2.6
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 6.5
$ws.Range("P4").Value = 2.37
$ws.Range("Q4").Value = 2.7
$ws.Range("R4").Value = 1.44
$ws.Range("S4").Value = 1.54
$ws.Range("AB5").Value = 34
$ws.Range("AK5").Value = 21
$ws.Range("AR5").Value = 81
$ws.Range("AW5").Value = 4.33
$ws.Range("AY5").Value = 21
$ws.Range("BA5").Value = 51
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 2.3
$ws.Range("J5").Value = 3.6
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.27
$ws.Range("S5").Value = 1.37
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.95
$ws.Range("W5").Value = 9.5
$ws.Range("Z5").Value = 34
$ws.Range("AC6").Value = 13.5
$ws.Range("AE6").Value = 16
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 18
$ws.Range("AJ6").Value = 19.5
$ws.Range("AK6").Value = 150
$ws.Range("AP6").Value = 15.5
$ws.Range("AT6").Value = 3.2
$ws.Range("AV6").Value = 70
$ws.Range("AY6").Value = 40
$ws.Range("BB6").Value = 500
$ws.Range("G6").Value = 1.32
$ws.Range("H6").Value = 4.55
$ws.Range("I6").Value = 7.5
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 12.8
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 4.15
$ws.Range("T6").Value = 3.32
$ws.Range("W6").Value = 6.4
$ws.Range("X6").Value = 5.7
$ws.Range("Z6").Value = 7.2
$ws.Range("AA7").Value = 27
$ws.Range("AB7").Value = 37
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 5.4
$ws.Range("AE7").Value = 14
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 600
$ws.Range("AH7").Value = 7
$ws.Range("AI7").Value = 12
$ws.Range("AJ7").Value = 9.75
$ws.Range("AK7").Value = 30
$ws.Range("AL7").Value = 24
$ws.Range("AM7").Value = 35
$ws.Range("AN7").Value = 4.8
$ws.Range("AO7").Value = 16
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 75
$ws.Range("AR7").Value = 100
$ws.Range("AT7").Value = 2.45
$ws.Range("AU7").Value = 6.4
$ws.Range("AW7").Value = 4.45
$ws.Range("AX7").Value = 14
$ws.Range("G7").Value = 2.95
$ws.Range("H7").Value = 2.77
$ws.Range("I7").Value = 2.57
$ws.Range("J7").Value = 3.4
$ws.Range("K7").Value = 1.98
$ws.Range("L7").Value = 3.15
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 6.8
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.55
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.55
$ws.Range("T7").Value = 2.47
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.82
$ws.Range("W7").Value = 8
$ws.Range("X7").Value = 15
$ws.Range("Y7").Value = 10.5
$ws.Range("U8").Value = 1.69
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 15
$ws.Range("AD9").Value = 8
$ws.Range("AH9").Value = 19
$ws.Range("AI9").Value = 34
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 67
$ws.Range("BD9").Value = 176
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.7
$ws.Range("R9").Value = 2.1
$ws.Range("U9").Value = 1.69
$ws.Range("W9").Value = 8.5
$ws.Range("X9").Value = 8
$ws.Range("Z9").Value = 11
$ws.Range("AB10").Value = 29
$ws.Range("AH10").Value = 10
$ws.Range("AK10").Value = 34
$ws.Range("AL10").Value = 26
$ws.Range("AN10").Value = 4.33
$ws.Range("AO10").Value = 12
$ws.Range("AW10").Value = 5
$ws.Range("AY10").Value = 26
$ws.Range("AZ10").Value = 51
$ws.Range("BC10").Value = 126
$ws.Range("G10").Value = 2.15
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 2.88
$ws.Range("K10").Value = 2.1
$ws.Range("L10").Value = 3.75
$ws.Range("W10").Value = 8
$ws.Range("Z10").Value = 21
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("AA13").Value = 13
$ws.Range("AB13").Value = 29
$ws.Range("AH13").Value = 21
$ws.Range("AN13").Value = 3.25
$ws.Range("G13").Value = 1.38
$ws.Range("H13").Value = 4.2
$ws.Range("I13").Value = 9
$ws.Range("J13").Value = 1.91
$ws.Range("K13").Value = 2.38
$ws.Range("L13").Value = 7.5
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.95
$ws.Range("R13").Value = 1.9
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
